# Edit PO Forecast: rename headers on existing sheets, add new "PO Forecast" sheet.

$wb = $excel.ActiveWorkbook

# 1) Rename "Requested quantity" header on "Weekly Quantity" -> "Weekly_PO_Qty"
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2) Rename "Requested quantity" header on "Monthly Trend" -> "Monthly_PO_Qty"
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3) Add new "PO Forecast" worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Cells.Item(1,1).Value = "ds"
$wsForecast.Cells.Item(1,2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1,3).Value = "yhat_lower"
$wsForecast.Cells.Item(1,4).Value = "yhat_upper"

# Match header formatting (bold / centered / bordered) used on the other sheets.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$rows = @(
    @(44948.99999999999,72,-23.17870295691156,171.2464854563795),
    @(44969.99999999999,73,-22.31227442871758,175.8274711056413),
    @(44976.99999999999,73,-27.7933072492904,171.6393807630572),
    @(44983.99999999999,73,-27.57592564660463,178.6028190745696),
    @(44990.99999999999,74,-19.40546239322103,162.7451755532419),
    @(44997.99999999999,74,-19.04756806199492,171.5964373188206),
    @(45004.99999999999,74,-22.04123541458574,168.7068401717488),
    @(45011.99999999999,74,-17.18307259599341,170.6845228871883),
    @(45025.99999999999,75,-18.64649781873675,173.2109355932137),
    @(45053.99999999999,76,-28.02153132516399,179.9467448169417),
    @(45067.99999999999,77,-20.68711174533389,173.88550339502),
    @(45074.99999999999,77,-16.12256572095369,175.0201481840314),
    @(45081.99999999999,77,-17.25789275544621,174.4990766053326),
    @(45102.99999999999,78,-17.5947273267934,181.4724085642404),
    @(45109.99999999999,79,-22.75828070123155,168.0360761961781),
    @(45116.99999999999,79,-18.15390364201326,180.1492408351156),
    @(45137.99999999999,80,-14.65769044632731,172.0249109633792),
    @(45165.99999999999,81,-22.25244819240448,177.5244955863363),
    @(45172.99999999999,81,-12.0883199880502,170.8314747905214),
    @(45179.99999999999,82,-18.02878652998077,179.001145601558),
    @(45186.99999999999,82,-13.62420835100587,182.5244036955334),
    @(45193.99999999999,82,-9.86026198067955,179.6079001319349),
    @(45200.99999999999,82,-10.66786330326786,182.6896777926376),
    @(45207.99999999999,83,-11.29537399916927,182.985512910049),
    @(45214.99999999999,83,-14.77726267870409,188.5462845448174),
    @(45235.99999999999,84,-11.45743323474555,182.0883482547872),
    @(45242.99999999999,84,-9.877895180829263,184.2164482341189),
    @(45249.99999999999,84,-14.54517429320668,184.1706423782925),
    @(45256.99999999999,85,-17.36948094533476,178.052651059761),
    @(45263.99999999999,85,-9.971083642442894,185.0903662648958),
    @(45277.99999999999,86,-14.64723260546183,178.9532853405661),
    @(45298.99999999999,87,-9.137830801508045,186.6401356488127),
    @(45305.99999999999,87,-12.72335603536922,179.6793678067246),
    @(45319.99999999999,87,-16.86552923411721,186.1901388651395),
    @(45326.99999999999,88,-6.785253642049077,182.427768681053),
    @(45333.99999999999,88,-9.083269181425313,189.2323542926313),
    @(45340.99999999999,88,-16.64793980103079,188.9118081957498),
    @(45347.99999999999,89,-5.11325393535194,184.1268833497714),
    @(45354.99999999999,89,-8.875248256846,183.4190617283847),
    @(45361.99999999999,89,-4.377847138876242,178.8837436244571),
    @(45389.99999999999,90,-11.13010655224764,188.8231036114998),
    @(45396.99999999999,91,-5.927133581976693,187.2002514624511),
    @(45403.99999999999,91,-6.892789592340905,187.7806443978303),
    @(45410.99999999999,91,-12.40695061628643,185.8839003828467),
    @(45417.99999999999,92,-13.08103676298393,187.7698916984506),
    @(45424.99999999999,92,-9.746809382194606,184.7228682017362),
    @(45431.99999999999,92,-1.824874586513942,190.2207435045006),
    @(45445.99999999999,93,-6.115776034537381,198.1478083279539),
    @(45459.99999999999,93,0.9484728811974799,196.339188377058),
    @(45487.99999999999,94,0.09812852961886219,191.8309302526625),
    @(45515.99999999999,96,-2.601873739396416,187.2668146662218),
    @(45522.99999999999,96,-6.595644199934619,191.3181314624053),
    @(45529.99999999999,96,-3.670834849154157,189.0977748156473),
    @(45536.99999999999,97,0.02111201053453136,189.1183521752206),
    @(45543.99999999999,97,1.237007718251721,193.8361992595586),
    @(45550.99999999999,97,-7.593477681770353,188.1977140463177),
    @(45557.99999999999,97,2.849032166795295,200.439943765528),
    @(45564.99999999999,98,4.716654214327624,202.810071239922),
    @(45571.99999999999,98,-1.460514827307616,198.2297435506051),
    @(45578.99999999999,98,2.641610318890103,196.3869815240418)
)

$n = $rows.Count
$data = New-Object 'object[,]' $n,4
for ($i = 0; $i -lt $n; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$wsForecast.Range("A2:D61").Value = $data

# Match the date-column formatting (yyyy-mm-dd hh:mm:ss number format) used by
# the "ds"/"Order Week" date columns on the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A61").PasteSpecial(-4122)

Write-Output "PO Forecast sheet added; headers renamed."
